$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add a new "Total" column to the Stock table (extends A1:C4 -> A1:D4)
$lo.ListColumns.Add()
$ws.Range("D1").Value = "Total"

# Replace the structured-reference formulas with plain cell-reference
# formulas (fixes corrupted table formulas)
$ws.Range("D2").Formula = "=B2*C2"
$ws.Range("D3").Formula = "=B3*C3"
$ws.Range("D4").Formula = "=B4*C4"

# Restore table style: show first column banding, no row stripes
$lo.ShowTableStyleFirstColumn = $true
$lo.ShowTableStyleRowStripes = $false
